$d = $word.ActiveDocument

# --- Replace the title paragraph + the "By Dorothy Day" paragraph ---
# Original:
#   P1 (Heading1, bookmarked): "The Daily Worker Case"
#   P2 (Normal, bold run):     "By Dorothy Day"
# Target:
#   P1 (Title): "The" " " "Daily" " " "Worker" " " "Case"   (split into separate runs)
#   P2 (Authors): "Dorothy" " " "Day"                        (split into separate runs)

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$titleRuns = @("The", " ", "Daily", " ", "Worker", " ", "Case")
$titleRunsXml = ""
foreach ($t in $titleRuns) {
    $titleRunsXml += "<w:r><w:t xml:space=`"preserve`">$t</w:t></w:r>"
}

$authorRuns = @("Dorothy", " ", "Day")
$authorRunsXml = ""
foreach ($t in $authorRuns) {
    $authorRunsXml += "<w:r><w:t xml:space=`"preserve`">$t</w:t></w:r>"
}

$newXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$titleRunsXml</w:p>" +
          "<w:p $wNs><w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>$authorRunsXml</w:p>"

$full = $d.Range($p1.Range.Start, $p2.Range.End)
$full.InsertXML($newXml)

# --- Best-effort removal of the stray "the-daily-worker-case" bookmark ---
# (The source document wraps the old title paragraph in a bookmark; the
# target no longer has it.)
if ($d.Bookmarks.Exists("the-daily-worker-case")) {
    $d.Bookmarks("the-daily-worker-case").Delete()
}
$d.DeleteBookmark("the-daily-worker-case")
